$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.603.84"
$ws.Range("E2").Value = "  +2.27%  "
$ws.Range("D3").Value = "3.833.27"
$ws.Range("E3").Value = "  +1.10%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.42%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "634.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.25"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.51%  "
$ws.Range("D7").Value = "3.831.25"
$ws.Range("E7").Value = "  +1.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +0.65%  "
$ws.Range("E10").Value = "  +1.98%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.454"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.59%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.68"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.12%  "
$ws.Range("E13").Value = "  +0.93%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.00"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.64%  "
$ws.Range("D15").Value = "4.475.82"
$ws.Range("E15").Value = "  +1.13%  "
$ws.Range("D16").Value = "3.923.29"
$ws.Range("E16").Value = "  +3.82%  "
$ws.Range("D17").Value = "69.562.13"
$ws.Range("E17").Value = "  +2.25%  "
$ws.Range("E18").Value = "  -1.22%  "
$ws.Range("E19").Value = "  +1.21%  "
$ws.Range("E20").Value = "  -0.29%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "466.98"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.71"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.709"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.66%  "
$ws.Range("E24").Value = "  +2.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.61"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.99%  "
$ws.Range("E26").Value = "  +3.53%  "
$ws.Range("E27").Value = "  -0.29%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.93%  "
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("D30").Value = "3.983.14"
$ws.Range("E30").Value = "  +1.09%  "
$ws.Range("E31").Value = "  +2.80%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.31"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.11%  "
$ws.Range("E33").Value = "  -0.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.27"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.19%  "
$ws.Range("D35").Value = "3.777.14"
$ws.Range("E35").Value = "  +0.97%  "
$ws.Range("E36").Value = "  +0.11%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "9.07"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.75%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.103"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.150"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +8.03%  "
$ws.Range("E40").Value = "  +6.93%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.92"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.53%  "
$ws.Range("E42").Value = "  -0.87%  "
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "157.47"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.31%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "44.09"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.54%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.302"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.82%  "
$ws.Range("E48").Value = "  +5.12%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "47.20"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.32%  "
$ws.Range("E50").Value = "  +3.11%  "
$ws.Range("E51").Value = "  +1.23%  "
